# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.139.86"
$ws.Range("E2").Value = "  -3.20%  "

# Row 3
$ws.Range("D3").Value = "1.711.69"
$ws.Range("E3").Value = "  -3.62%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.90%  "

# Row 6
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4768"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.88%  "

# Row 8
$ws.Range("E8").Value = "  -3.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07282"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.99%  "

# Row 12
$ws.Range("E12").Value = "  +0.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.853"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "

# Row 15
$ws.Range("D15").Value = "1.711.10"
$ws.Range("E15").Value = "  -3.50%  "

# Row 16
$ws.Range("E16").Value = "  -5.93%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.84%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001042"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06361"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "

# Row 20
$ws.Range("E20").Value = "  +0.16%  "

# Row 21
$ws.Range("E21").Value = "  -3.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.616"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.00%  "

# Row 23
$ws.Range("D23").Value = "27.177.52"
$ws.Range("E23").Value = "  -3.05%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.50%  "

# Row 25
$ws.Range("E25").Value = "  -2.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.82%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.78%  "

# Row 28
$ws.Range("D28").Value = "1.907.87"
$ws.Range("E28").Value = "  -3.50%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.086"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.014"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.81%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09260"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.589"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.299"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02197"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05896"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.15%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2014"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.754"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.411"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.92%  "

# Row 41
$ws.Range("E41").Value = "  +0.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5925"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.29%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.112"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.34%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.504"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.18%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.76%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.569"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5620"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.53%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.840"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.31%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06634"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.087"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.98%  "
